# Update "想去人数" (interest count, column F) for several events across
# the 展览 (Exhibition), 演出 (Performance) and 全部类型 (All types) sheets.
# Mirrors the regenerated data output (gh-pages build at 456a3b4).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 2795   # 北京·thebONE✖️GOJO超次元嘉年华12nd
$ws1.Range("F5").Value  = 6373   # 北京·第22届中国国际模型博览会
$ws1.Range("F13").Value = 6976   # 北京·IDO动漫游戏嘉年华45th
$ws1.Range("F19").Value = 8091   # 北京·第16届IJOY漫展XCGF游戏节
$ws1.Range("F41").Value = 166    # 北京·ICOS国际动漫节×CGF中国游戏节02
$ws1.Range("F43").Value = 145    # 北京·IDO动漫游戏嘉年华46th

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value  = 35     # 北京·acg乐队演出·你不会是红白歌会吧

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 35     # 北京·acg乐队演出·你不会是红白歌会吧
$ws4.Range("F3").Value  = 2795   # 北京·thebONE✖️GOJO超次元嘉年华12nd
$ws4.Range("F7").Value  = 6373   # 北京·第22届中国国际模型博览会
$ws4.Range("F18").Value = 6976   # 北京·IDO动漫游戏嘉年华45th
$ws4.Range("F24").Value = 8091   # 北京·第16届IJOY漫展XCGF游戏节
$ws4.Range("F47").Value = 166    # 北京·ICOS国际动漫节×CGF中国游戏节02
$ws4.Range("F50").Value = 145    # 北京·IDO动漫游戏嘉年华46th

$wb.Save()
